# Apply lattice multiplication exercise updates: replace the 15 cell
# contents (problem + partial products) in the 5x3 table in place,
# using $vtab (Word's internal manual-line-break char, matching <w:br/>)
# to join each cell's 5 lines so the whole cell is set atomically
# (avoids ambiguous Find/Replace matches when a cell has two identical
# lines, e.g. row 3 col 1 which originally has "2|    |" twice).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vtab = [char]11

# Row 1 Col 1: "80 x 29" -> "17 x 44"
$t.Cell(1, 1).Range.Text = "17 x 44" + $vtab + "  4    4" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "7|    |"

# Row 1 Col 2: "78 x 92" -> "58 x 33"
$t.Cell(1, 2).Range.Text = "58 x 33" + $vtab + "  3    3" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "8|    |"

# Row 1 Col 3: "18 x 82" -> "52 x 95"
$t.Cell(1, 3).Range.Text = "52 x 95" + $vtab + "  9    5" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "2|    |"

# Row 2 Col 1: "35 x 46" -> "10 x 93"
$t.Cell(2, 1).Range.Text = "10 x 93" + $vtab + "  9    3" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "0|    |"

# Row 2 Col 2: "15 x 23" -> "30 x 91"
$t.Cell(2, 2).Range.Text = "30 x 91" + $vtab + "  9    1" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "0|    |"

# Row 2 Col 3: "71 x 15" -> "82 x 53"
$t.Cell(2, 3).Range.Text = "82 x 53" + $vtab + "  5    3" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "2|    |"

# Row 3 Col 1: "22 x 11" -> "19 x 82"
$t.Cell(3, 1).Range.Text = "19 x 82" + $vtab + "  8    2" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "9|    |"

# Row 3 Col 2: "69 x 59" -> "78 x 68"
$t.Cell(3, 2).Range.Text = "78 x 68" + $vtab + "  6    8" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "8|    |"

# Row 3 Col 3: "51 x 34" -> "50 x 54"
$t.Cell(3, 3).Range.Text = "50 x 54" + $vtab + "  5    4" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "0|    |"

# Row 4 Col 1: "19 x 52" -> "25 x 32"
$t.Cell(4, 1).Range.Text = "25 x 32" + $vtab + "  3    2" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "5|    |"

# Row 4 Col 2: "48 x 29" -> "64 x 60"
$t.Cell(4, 2).Range.Text = "64 x 60" + $vtab + "  6    0" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "4|    |"

# Row 4 Col 3: "90 x 56" -> "22 x 10"
$t.Cell(4, 3).Range.Text = "22 x 10" + $vtab + "  1    0" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "2|    |"

# Row 5 Col 1: "10 x 39" -> "62 x 83"
$t.Cell(5, 1).Range.Text = "62 x 83" + $vtab + "  8    3" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "2|    |"

# Row 5 Col 2: "96 x 32" -> "94 x 21"
$t.Cell(5, 2).Range.Text = "94 x 21" + $vtab + "  2    1" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "4|    |"

# Row 5 Col 3: "83 x 51" -> "84 x 81"
$t.Cell(5, 3).Range.Text = "84 x 81" + $vtab + "  8    1" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "4|    |"
